# correction in sa algorithm and 746 logs
# Update the "Fitness" column (C) for generations 0..142 (rows 2..144) of run 29
# to reflect the corrected simulated-annealing algorithm output.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C21").Value = 7883
$ws.Range("C22:C26").Value = 7815
$ws.Range("C27:C49").Value = 7767
$ws.Range("C50:C53").Value = 7765
$ws.Range("C54:C55").Value = 7318
$ws.Range("C56:C144").Value = 7293
